# section_test.xlsx - "settings" sheet:
# Insert two new setting rows ("section1" / "section2") right after the
# "survey" row (row 4) and before the "default"/"hindi" rows, so the
# pseudo-prompt for each section gets a title shown in the form contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Push the existing "default" (row 5) / "hindi" (row 6) rows down to make
# room for the two new rows.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# setting_name column
$ws.Range("A5").Value() = "section1"
$ws.Range("A6").Value() = "section2"

# display.title / display.title.hindi columns
$ws.Range("C6").Value() = "Section 2"
$ws.Range("C5").Value() = "Section 1"

$ws.Range("D5").Value() = "धारा 1"
$ws.Range("D6").Value() = "धारा 2"

# Move the active selection the way it ended up after the edit.
$ws.Range("C10").Select() | Out-Null
